$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.055285096168518
$ws.Range("B1").Value = 3.332328319549561
$ws.Range("C1").Value = 2.876245021820068
$ws.Range("D1").Value = 2.297162532806396
$ws.Range("E1").Value = 1.46152675151825
